$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phase 1: preserve formatting of the two existing data rows (old row3 =
# "Doc Retrieval" / SIPL5316 pattern, old row4 = "Typing" pattern) into their
# new destination rows (6 and 9) BEFORE their content gets overwritten by the
# re-inserted "Search & Typing" rows.
# ---------------------------------------------------------------------------
$ws.Range("A3:N3").Copy()
$ws.Range("A6:N6").PasteSpecial(-4122)

$ws.Range("A4:N4").Copy()
$ws.Range("A9:N9").PasteSpecial(-4122)

# Rows 3,4,5,7,8 become new "Search & Typing" rows matching row 2's pattern,
# so clone row 2's formatting into them.
$ws.Range("A2:N2").Copy()
$ws.Range("A3:N3").PasteSpecial(-4122)
$ws.Range("A4:N4").PasteSpecial(-4122)
$ws.Range("A5:N5").PasteSpecial(-4122)
$ws.Range("A7:N7").PasteSpecial(-4122)
$ws.Range("A8:N8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Phase 2: write cell values row-by-row (row order matters: it determines the
# order new shared strings are appended in, which must mirror the original
# table this revert restores).
# ---------------------------------------------------------------------------

# Row 3
$ws.Range("A3").Value = 45437.041666608799
$ws.Range("B3").Value = "ORDS18-002"
$ws.Range("C3").Value = "SIPL6118"
$ws.Range("D3").Value = "SIPL4167"
$ws.Range("E3").Value = "SIPL5317"
$ws.Range("F3").Value = "SIPL5317"
$ws.Range("G3").Value = "Old Republic Diversified Services"
$ws.Range("H3").Value = "Title"
$ws.Range("I3").Value = "Search & Typing"
$ws.Range("J3").Value = "Update Search"
$ws.Range("K3").Value = "AL"
$ws.Range("L3").Value = "Shelby"
$ws.Range("M3").Value = "WIP"
$ws.Range("N3").Clear()
$ws.Range("N3").Value = ""

# Row 4
$ws.Range("A4").Value = 45438.041666608799
$ws.Range("B4").Value = "ORDS18-003"
$ws.Range("C4").Value = "SIPL6118"
$ws.Range("D4").Value = "SIPL4167"
$ws.Range("E4").Value = "SIPL5317"
$ws.Range("F4").Value = "SIPL5317"
$ws.Range("G4").Value = "Old Republic Diversified Services"
$ws.Range("H4").Value = "Title"
$ws.Range("I4").Value = "Search & Typing"
$ws.Range("J4").Value = "30 Years Search / Full Search"
$ws.Range("K4").Value = "AL"
$ws.Range("L4").Value = "Shelby"
$ws.Range("M4").Value = "WIP"
$ws.Range("N4").Clear()
$ws.Range("N4").Value = ""

# Row 5
$ws.Range("A5").Value = 45439.041666608799
$ws.Range("B5").Value = "ORDS18-004"
$ws.Range("C5").Value = "SIPL6118"
$ws.Range("D5").Value = "SIPL4167"
$ws.Range("E5").Value = "SIPL5317"
$ws.Range("F5").Value = "SIPL5317"
$ws.Range("G5").Value = "Old Republic Diversified Services"
$ws.Range("H5").Value = "Title"
$ws.Range("I5").Value = "Search & Typing"
$ws.Range("J5").Value = "Old Republic Diversified Services - ABS"
$ws.Range("K5").Value = "AL"
$ws.Range("L5").Value = "Shelby"
$ws.Range("M5").Value = "WIP"
$ws.Range("N5").Clear()
$ws.Range("N5").Value = ""

# Row 6 (restores the pre-existing "Doc Retrieval" row, shifted down)
$ws.Range("A6").Value = 45440.041666608799
$ws.Range("B6").Value = "ORDS18-005"
$ws.Range("C6").Value = "SIPL5316"
$ws.Range("D6").Value = "SIPL5688"
$ws.Range("G6").Value = "Old Republic Diversified Services"
$ws.Range("H6").Value = "Title"
$ws.Range("I6").Value = "Search"
$ws.Range("J6").Value = "Doc Retrieval"
$ws.Range("K6").Value = "AL"
$ws.Range("L6").Value = "Autauga"
$ws.Range("M6").Value = "WIP"
$ws.Range("N6").Value = "Search(T1)"

# Row 7
$ws.Range("A7").Value = 45441.041666608799
$ws.Range("B7").Value = "ORDS18-006"
$ws.Range("C7").Value = "SIPL6118"
$ws.Range("D7").Value = "SIPL4167"
$ws.Range("E7").Value = "SIPL5317"
$ws.Range("F7").Value = "SIPL5317"
$ws.Range("G7").Value = "Old Republic Diversified Services"
$ws.Range("H7").Value = "Title"
$ws.Range("I7").Value = "Search & Typing"
$ws.Range("J7").Value = "L&V with Taxes"
$ws.Range("K7").Value = "AL"
$ws.Range("L7").Value = "Shelby"
$ws.Range("M7").Value = "WIP"
$ws.Range("N7").Clear()
$ws.Range("N7").Value = ""

# Row 8
$ws.Range("A8").Value = 45442.041666608799
$ws.Range("B8").Value = "ORDS18-007"
$ws.Range("C8").Value = "SIPL6118"
$ws.Range("D8").Value = "SIPL4167"
$ws.Range("E8").Value = "SIPL5317"
$ws.Range("F8").Value = "SIPL5317"
$ws.Range("G8").Value = "Old Republic Diversified Services"
$ws.Range("H8").Value = "Title"
$ws.Range("I8").Value = "Search & Typing"
$ws.Range("J8").Value = "COS - Project"
$ws.Range("K8").Value = "AL"
$ws.Range("L8").Value = "Shelby"
$ws.Range("M8").Value = "WIP"
$ws.Range("N8").Clear()
$ws.Range("N8").Value = ""

# Row 9 (restores the pre-existing "Typing" row, shifted down)
$ws.Range("A9").Value = 45443.041666608799
$ws.Range("B9").Value = "ORDS18-008"
$ws.Range("C9").Value = "SIPL5316"
$ws.Range("D9").Value = "SIPL5688"
$ws.Range("E9").Clear()
$ws.Range("F9").Clear()
$ws.Range("G9").Value = "Old Republic Diversified Services"
$ws.Range("H9").Value = "Title"
$ws.Range("I9").Value = "Typing"
$ws.Range("J9").Value = "Typing"
$ws.Range("K9").Value = "AL"
$ws.Range("L9").Value = "Autauga"
$ws.Range("M9").Value = "WIP"
$ws.Range("N9").Value = "Typing(T1)"

Write-Host "values written"
